# Auto-generated Excel COM-interop script to apply the Gilgamesh_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1550
$ws.Cells.Item(2, 10).Value = 1583.3334
$ws.Cells.Item(2, 12).Value = 1583.3334
$ws.Cells.Item(2, 14).Value = -1809.3334
$ws.Cells.Item(87, 8).Value = 254999
$ws.Cells.Item(87, 10).Value = 254999
$ws.Cells.Item(87, 12).Value = 254999
$ws.Cells.Item(87, 14).Value = -257495
$ws.Cells.Item(90, 8).Value = 254999
$ws.Cells.Item(90, 10).Value = 254999
$ws.Cells.Item(90, 12).Value = 764997
$ws.Cells.Item(90, 14).Value = -777477
$ws.Cells.Item(92, 8).Value = 439.35
$ws.Cells.Item(92, 9).Value = 410.4375
$ws.Cells.Item(92, 10).Value = 555
$ws.Cells.Item(92, 11).Value = 410.4375
$ws.Cells.Item(92, 12).Value = 555
$ws.Cells.Item(92, 13).Value = 837.5625
$ws.Cells.Item(92, 14).Value = -3051
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 389.57144
$ws.Cells.Item(96, 9).Value = 358
$ws.Cells.Item(96, 11).Value = 1074
$ws.Cells.Item(96, 13).Value = 299
$ws.Cells.Item(97, 8).Value = 12254.8
$ws.Cells.Item(97, 10).Value = 12254.8
$ws.Cells.Item(97, 12).Value = 36764.39999999999
$ws.Cells.Item(97, 14).Value = -37756.39999999999
$ws.Cells.Item(107, 8).Value = 457.8
$ws.Cells.Item(107, 10).Value = 951
$ws.Cells.Item(107, 12).Value = 951
$ws.Cells.Item(107, 14).Value = -4791
$ws.Cells.Item(127, 8).Value = 491.33334
$ws.Cells.Item(127, 9).Value = 513.36365
$ws.Cells.Item(127, 11).Value = 1540.09095
$ws.Cells.Item(127, 13).Value = 3419.90905
$ws.Cells.Item(138, 8).Value = 4392.2173
$ws.Cells.Item(138, 10).Value = 3779.8809
$ws.Cells.Item(138, 12).Value = 11339.6427
$ws.Cells.Item(138, 14).Value = -21619.6427

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 1387.0769
$ws.Cells.Item(97, 9).Value = 1139.2727
$ws.Cells.Item(97, 10).Value = 2750
$ws.Cells.Item(97, 11).Value = 1139.2727
$ws.Cells.Item(97, 12).Value = 2750
$ws.Cells.Item(97, 13).Value = -643.2727
$ws.Cells.Item(97, 14).Value = -3742
$ws.Cells.Item(122, 8).Value = 2127.2188
$ws.Cells.Item(122, 9).Value = 2127.2188
$ws.Cells.Item(122, 11).Value = 6381.6564
$ws.Cells.Item(122, 13).Value = -3931.6564
$ws.Cells.Item(132, 8).Value = 3083.5435
$ws.Cells.Item(132, 9).Value = 2730.8235
$ws.Cells.Item(132, 10).Value = 4082.9167
$ws.Cells.Item(132, 11).Value = 8192.470499999999
$ws.Cells.Item(132, 12).Value = 12248.7501
$ws.Cells.Item(132, 13).Value = -5662.470499999999
$ws.Cells.Item(132, 14).Value = -17308.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 105264696
$ws.Cells.Item(94, 9).Value = 166667420
$ws.Cells.Item(94, 10).Value = 2887.1428
$ws.Cells.Item(94, 11).Value = 166667420
$ws.Cells.Item(94, 12).Value = 2887.1428
$ws.Cells.Item(94, 13).Value = -166666969
$ws.Cells.Item(94, 14).Value = -3789.1428
$ws.Cells.Item(95, 8).Value = 62481.668
$ws.Cells.Item(95, 10).Value = 62481.668
$ws.Cells.Item(95, 12).Value = 62481.668
$ws.Cells.Item(95, 14).Value = -67973.66800000001
$ws.Cells.Item(105, 8).Value = 10835621
$ws.Cells.Item(105, 9).Value = 668360.9
$ws.Cells.Item(105, 10).Value = 27781054
$ws.Cells.Item(105, 11).Value = 668360.9
$ws.Cells.Item(105, 12).Value = 27781054
$ws.Cells.Item(105, 13).Value = -666613.9
$ws.Cells.Item(105, 14).Value = -27784548

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 68090
$ws.Cells.Item(50, 10).Value = 68090
$ws.Cells.Item(50, 12).Value = 68090
$ws.Cells.Item(50, 14).Value = -69340
$ws.Cells.Item(59, 8).Value = 88126.5
$ws.Cells.Item(59, 10).Value = 88126.5
$ws.Cells.Item(59, 12).Value = 88126.5
$ws.Cells.Item(59, 14).Value = -90416.5
$ws.Cells.Item(132, 8).Value = 4069.3103
$ws.Cells.Item(132, 9).Value = 3546.9443
$ws.Cells.Item(132, 11).Value = 10640.8329
$ws.Cells.Item(132, 13).Value = -8110.832900000001
$ws.Cells.Item(141, 8).Value = 476270.62
$ws.Cells.Item(141, 10).Value = 476270.62
$ws.Cells.Item(141, 12).Value = 476270.62
$ws.Cells.Item(141, 14).Value = -486630.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 2899.5
$ws.Cells.Item(92, 9).Value = 2849.75
$ws.Cells.Item(92, 10).Value = 2949.25
$ws.Cells.Item(92, 11).Value = 8549.25
$ws.Cells.Item(92, 12).Value = 8847.75
$ws.Cells.Item(92, 13).Value = -7301.25
$ws.Cells.Item(92, 14).Value = -11343.75
$ws.Cells.Item(93, 8).Value = 2560
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 13).ClearContents()
$ws.Cells.Item(97, 8).Value = 1250312.2
$ws.Cells.Item(97, 9).Value = 1250312.2
$ws.Cells.Item(97, 11).Value = 3750936.6
$ws.Cells.Item(97, 13).Value = -3750440.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2792.077
$ws.Cells.Item(97, 9).Value = 2299.6667
$ws.Cells.Item(97, 10).Value = 3900
$ws.Cells.Item(97, 11).Value = 2299.6667
$ws.Cells.Item(97, 12).Value = 3900
$ws.Cells.Item(97, 13).Value = -1803.6667
$ws.Cells.Item(97, 14).Value = -4892
$ws.Cells.Item(113, 8).Value = 1974.05
$ws.Cells.Item(113, 9).Value = 1835.8
$ws.Cells.Item(113, 11).Value = 1835.8
$ws.Cells.Item(113, 13).Value = 334.2
$ws.Cells.Item(122, 8).Value = 3555.4243
$ws.Cells.Item(122, 9).Value = 2574.96
$ws.Cells.Item(122, 10).Value = 6619.375
$ws.Cells.Item(122, 11).Value = 7724.88
$ws.Cells.Item(122, 12).Value = 19858.125
$ws.Cells.Item(122, 13).Value = -5274.88
$ws.Cells.Item(122, 14).Value = -24758.125
$ws.Cells.Item(132, 8).Value = 3219.182
$ws.Cells.Item(132, 9).Value = 2344.6428
$ws.Cells.Item(132, 10).Value = 4749.625
$ws.Cells.Item(132, 11).Value = 7033.928400000001
$ws.Cells.Item(132, 12).Value = 14248.875
$ws.Cells.Item(132, 13).Value = -4503.928400000001
$ws.Cells.Item(132, 14).Value = -19308.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 1359.4
$ws.Cells.Item(9, 9).Value = 199.25
$ws.Cells.Item(9, 11).Value = 199.25
$ws.Cells.Item(9, 13).Value = 24.75
$ws.Cells.Item(16, 8).Value = 852
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 2457.4
$ws.Cells.Item(46, 9).Value = 2522
$ws.Cells.Item(46, 11).Value = 2522
$ws.Cells.Item(46, 13).Value = -2334
$ws.Cells.Item(93, 8).Value = 550
$ws.Cells.Item(93, 9).Value = 550
$ws.Cells.Item(93, 11).Value = 550
$ws.Cells.Item(93, 13).Value = 698
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 1421.75
$ws.Cells.Item(23, 9).Value = 1662.6666
$ws.Cells.Item(23, 11).Value = 1662.6666
$ws.Cells.Item(23, 13).Value = -1433.6666
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 3878.6
$ws.Cells.Item(96, 9).Value = 3464.6667
$ws.Cells.Item(96, 11).Value = 3464.6667
$ws.Cells.Item(96, 13).Value = -2091.6667
$ws.Cells.Item(103, 8).Value = 92367.664
$ws.Cells.Item(103, 10).Value = 92367.664
$ws.Cells.Item(103, 12).Value = 92367.664
$ws.Cells.Item(103, 14).Value = -94711.664
$ws.Cells.Item(113, 8).Value = 1315.5358
$ws.Cells.Item(113, 9).Value = 1247.1875
$ws.Cells.Item(113, 10).Value = 1406.6666
$ws.Cells.Item(113, 11).Value = 3741.5625
$ws.Cells.Item(113, 12).Value = 4219.9998
$ws.Cells.Item(113, 13).Value = -1571.5625
$ws.Cells.Item(113, 14).Value = -8559.9998
$ws.Cells.Item(126, 8).Value = 2987.2222
$ws.Cells.Item(126, 9).Value = 3020.8462
$ws.Cells.Item(126, 11).Value = 9062.5386
$ws.Cells.Item(126, 13).Value = -6592.5386
